$wb = $excel.ActiveWorkbook
$ws4 = $wb.Worksheets.Item("Week4")
$ws4.Copy($null, $ws4)
$new = $wb.Worksheets.Item($wb.Worksheets.Count)
$new.Name = "Week5"

# Update dates row 5 B:H -> next week
$new.Range("B5").Value = 45221
$new.Range("C5").Value = 45222
$new.Range("D5").Value = 45223
$new.Range("E5").Value = 45224
$new.Range("F5").Value = 45225
$new.Range("G5").Value = 45226
$new.Range("H5").Value = 45227

# Update hours for new week
$new.Range("E8").Value = 1
$new.Range("E10").Value = $null
$new.Range("D11").Value = $null
$new.Range("E11").Value = $null
$new.Range("F11").Value = 2
$new.Range("G12").Value = $null

# Note
$new.Range("I2").Value = "*3 midterms this week, less time than usual"

# Remove spare column J (Week5 matches Week1-3 layout with no J column)
$new.Range("J1:J24").EntireColumn.Delete()

Write-Host "done basic edits"
